$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old "Name"/"Enrollment" columns are gone; the "12 aug" attendance
# column becomes the sole (new) column A. Deleting columns A:B shifts the
# old column C ("12 aug", P/absent marks) left into column A, which is the
# dataset that was "added" (and trained on).
$ws.Columns("A:B").Delete()

# Row 3 ("rishi") had no attendance mark in the original sheet; the
# refreshed dataset records "P" for that day.
$ws.Range("A3").Value = "P"
